# Update the "path" data table on Sheet1: refresh existing rows' values and
# append new rows 5-7, then move the active selection to E8 to match the
# edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (id=1): path gets extended.
$ws.Range("B2").Value = "1,2,3,5,4,6"

# Row 3 (id=2): path unchanged in content.
$ws.Range("B3").Value = "2,3,5"

# New row 5 (id=4) is entered before row 4 is updated, matching the order
# the values were authored in.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "4,6"

# Row 4 (id=3): path updated.
$ws.Range("B4").Value = "4,5"

# New row 6 (id=5).
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "1,8,11"

# New row 7 (id=6).
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "1,9,12"

# Match the saved selection state from the edited workbook.
$ws.Range("E8").Select()
